# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed values, matching the upstream GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, exactly as the source XML has it
# (inline/shared string, no numeric coercion, no lingering custom number
# format left behind on the cell).
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue 'D2' '27.314.55'
Set-TextValue 'E2' '  -2.10%  '
Set-TextValue 'D3' '1.826.80'
Set-TextValue 'E3' '  -1.85%  '
Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  -1.18%  '
Set-TextValue 'D5' '314.63'
Set-TextValue 'E5' '  -2.17%  '
Set-TextValue 'E6' '  -1.15%  '
Set-TextValue 'D7' '0.4273'
Set-TextValue 'E7' '  -1.44%  '
Set-TextValue 'D8' '0.3702'
Set-TextValue 'E8' '  -2.64%  '
Set-TextValue 'D9' '0.07253'
Set-TextValue 'E9' '  -2.43%  '
Set-TextValue 'D10' '0.8641'
Set-TextValue 'E10' '  -2.60%  '
Set-TextValue 'D11' '21.07'
Set-TextValue 'E11' '  -3.13%  '
Set-TextValue 'D12' '1.827.08'
Set-TextValue 'E12' '  -1.95%  '
Set-TextValue 'D13' '6.712'
Set-TextValue 'E13' '  -1.02%  '
Set-TextValue 'D14' '0.07115'
Set-TextValue 'E14' '  -0.52%  '
Set-TextValue 'D15' '5.317'
Set-TextValue 'E15' '  -3.53%  '
Set-TextValue 'D16' '88.65'
Set-TextValue 'E16' '  +0.15%  '
Set-TextValue 'E17' '  -1.36%  '
Set-TextValue 'D18' '0.000008872'
Set-TextValue 'E18' '  -2.09%  '
Set-TextValue 'E19' '  -1.03%  '
Set-TextValue 'D20' '15.10'
Set-TextValue 'E20' '  -3.04%  '
Set-TextValue 'D21' '27.315.19'
Set-TextValue 'E21' '  -2.20%  '
Set-TextValue 'D22' '5.143'
Set-TextValue 'E22' '  -2.68%  '
Set-TextValue 'D23' '10.89'
Set-TextValue 'E23' '  -3.07%  '
Set-TextValue 'D24' '2.050.49'
Set-TextValue 'E24' '  -2.11%  '
Set-TextValue 'E25' '  -1.35%  '
Set-TextValue 'D26' '153.12'
Set-TextValue 'E26' '  -2.49%  '
Set-TextValue 'D27' '18.39'
Set-TextValue 'E27' '  -1.78%  '
Set-TextValue 'D28' '2.152'
Set-TextValue 'E28' '  +6.08%  '
Set-TextValue 'D29' '5.247'
Set-TextValue 'E29' '  -3.59%  '
Set-TextValue 'D30' '116.51'
Set-TextValue 'E30' '  -3.78%  '
Set-TextValue 'D31' '0.08904'
Set-TextValue 'E31' '  -1.09%  '
Set-TextValue 'D32' '1.201'
Set-TextValue 'E32' '  -3.35%  '
Set-TextValue 'D33' '0.7578'
Set-TextValue 'E33' '  -2.60%  '
Set-TextValue 'D34' '4.451'
Set-TextValue 'E34' '  -3.30%  '
Set-TextValue 'E35' '  -2.92%  '
Set-TextValue 'E36' '  -1.15%  '
Set-TextValue 'D37' '1.112'
Set-TextValue 'E37' '  -3.60%  '
Set-TextValue 'D38' '0.01974'
Set-TextValue 'E38' '  -0.15%  '
Set-TextValue 'D39' '0.05279'
Set-TextValue 'E39' '  -0.98%  '
Set-TextValue 'D40' '7.164'
Set-TextValue 'E40' '  +1.94%  '
Set-TextValue 'D41' '2.873'
Set-TextValue 'E41' '  -0.66%  '
Set-TextValue 'D42' '0.1697'
Set-TextValue 'E42' '  +0.66%  '
Set-TextValue 'D43' '0.5047'
Set-TextValue 'E43' '  -3.24%  '
Set-TextValue 'D44' '8.668'
Set-TextValue 'E44' '  -1.44%  '
Set-TextValue 'D45' '10.59'
Set-TextValue 'E45' '  -1.74%  '
Set-TextValue 'D46' '107.65'
Set-TextValue 'E46' '  -3.09%  '
Set-TextValue 'D47' '0.4756'
Set-TextValue 'E47' '  -0.25%  '
Set-TextValue 'E48' '  -1.25%  '
Set-TextValue 'D49' '0.06374'
Set-TextValue 'E49' '  -1.96%  '
Set-TextValue 'D50' '1.665'
Set-TextValue 'E50' '  -3.20%  '
Set-TextValue 'D51' '1.824'
Set-TextValue 'E51' '  -3.24%  '
